# Add team record (Wins/Losses/Ties) columns to the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: new columns AD, AE, AF
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the existing header style (bold, centered, bordered) used by A1:AC1
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Data rows 2-46: Wins=83, Losses=79, Ties=0 for every row
$lastRow = 46
$ws.Range("AD2:AD$lastRow").Value = 83
$ws.Range("AE2:AE$lastRow").Value = 79
$ws.Range("AF2:AF$lastRow").Value = 0
